$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows of data (2-6) have been re-sorted/re-ordered. The final values for
# each row below reflect that new order (Q/R coordinates are also rounded to
# whole numbers in the new version), the Starttid (Z) and Sluttid (AB) columns
# are cleared out entirely, and the "Publik kommentar" (AC) column now travels
# with its corresponding record.

# Row 2 (was id 111638283, now id 111638281)
$ws.Range("A2").Value = 111638281
$ws.Range("B2").Value = 89423
$ws.Range("E2").Value = 5432
$ws.Range("F2").Value = "Granticka"
$ws.Range("G2").Value = "Porodaedalea chrysoloma"
$ws.Range("H2").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q2").Value = 588278
$ws.Range("R2").Value = 7033317
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("AC2").ClearContents()

# Row 3 (was id 111638282, now id 111638277)
$ws.Range("A3").Value = 111638277
$ws.Range("B3").Value = 77267
$ws.Range("E3").Value = 6446
$ws.Range("F3").Value = "Kolflarnlav"
$ws.Range("G3").Value = "Carbonicola anthracophila"
$ws.Range("H3").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q3").Value = 588323
$ws.Range("R3").Value = 7033261
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
$ws.Range("AC3").Value = "På kolad tallstubbe med yxhugg."

# Row 4 (was id 111638281, now id 111638278)
$ws.Range("A4").Value = 111638278
$ws.Range("B4").Value = 56414
$ws.Range("E4").Value = 100049
$ws.Range("F4").Value = "Spillkråka"
$ws.Range("G4").Value = "Dryocopus martius"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("Q4").Value = 588320
$ws.Range("R4").Value = 7033286
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
$ws.Range("AC4").Value = "Bohål i gammal grov tall."

# Row 5 (was id 111638277, now id 111638283)
$ws.Range("A5").Value = 111638283
$ws.Range("B5").Value = 77515
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("Q5").Value = 588213
$ws.Range("R5").Value = 7033298
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()
$ws.Range("AC5").ClearContents()

# Row 6 (was id 111638278, now id 111638282)
$ws.Range("A6").Value = 111638282
$ws.Range("B6").Value = 89405
$ws.Range("E6").Value = 1202
$ws.Range("F6").Value = "Ullticka"
$ws.Range("G6").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H6").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q6").Value = 588266
$ws.Range("R6").Value = 7033312
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()
$ws.Range("AC6").ClearContents()
